# Update the LATAM CONSTRUCTION_STANDARD database: drop the unused
# "area_pv" and "area_sc" argument columns from ENVELOPE_ASSEMBLIES
# (commit message: "updating database to reduce arguments").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ENVELOPE_ASSEMBLIES")

# Columns W ("area_pv") and X ("area_sc") are removed entirely; the
# columns to their right (area_balcon, Refernez) shift left to take
# their place.
$ws.Range("W1:X1").EntireColumn.Delete()

# The sheet also had a stray formatted-but-empty row further down
# (only cell R9 carried a style, no real data) that goes away too.
$ws.Range("R9").EntireRow.Delete()
